$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.758.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.806.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5907"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.95%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06837"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07505"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.808.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.765"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6236"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.051.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009281"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "75.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.704.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.478"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.15%  "
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "211.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.26%  "
$ws.Range("E22").Value = "  -1.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.843"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.885"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1270"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.437"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06197"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.426"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.789"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.736"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.064"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6437"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.496"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.25%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.564"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01702"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.146.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8833"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.006"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.960.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.67%  "
$ws.Range("E47").Value = "  -5.42%  "
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.369"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05475"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4483"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.50%  "
